# Insert a new data row right after the existing "2026/01/30" entry (row 738),
# shifting the 2026/12/29 ... 2027/01/05 block down by one row (738:779 -> 739:780),
# and populate the newly-opened row with the extra 2026/01/30 reading (time 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 738:779 down to 739:780 by inserting a blank row at 738.
$ws.Rows.Item(738).Insert()

# Populate the newly inserted row 738.
# Column A holds date-looking text (e.g. "2026/01/30") stored as plain text in
# this sheet (same as every other row), so force text format before writing the
# value to avoid Excel auto-converting it to a date serial, then drop the
# temporary number format again so the cell keeps the sheet's default (no)
# style, just like its neighbours.
$ws.Cells.Item(738, 1).NumberFormat = "@"
$ws.Cells.Item(738, 1).Value = "2026/01/30"
$ws.Cells.Item(738, 1).Style = "Normal"

$ws.Cells.Item(738, 2).Value = "金"
$ws.Cells.Item(738, 3).Value = 8
$ws.Cells.Item(738, 4).Value = 201
